# Rename the Pearson / BTec logo pictures in the document's headers and
# footers so the inline picture's Name matches the new asset-numbering
# used by the authoring tool (image1.png <-> image2.png and
# image1.jpg <-> image2.jpg), without touching the alt-text/description,
# size, or any other property of the shapes.

$d = $word.ActiveDocument

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    # --- Headers: BTec logo (image2.jpg -> image1.jpg) ---
    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shp = $shapes.Item($j)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                }
            }
        }
    }

    # --- Footers: Pearson Edexcel logo (image1.png -> image2.png) ---
    for ($i = 1; $i -le 3; $i++) {
        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shp = $shapes.Item($j)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image2.png"
                }
            }
        }
    }
}
